$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update status text "In Translation" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Update timestamps to reflect the new handoff generation
$overview.Range("G2").Value = "2016-08-15 14:38:14"
$dede.Range("H2").Value = "2016-08-15 14:38:14"
$zhcn.Range("H2").Value = "2016-08-15 14:38:09"

# Widen columns to fit the new, longer status text
# (ColumnWidth is specified in characters and Excel snaps the stored width
# to whole-pixel increments at the workbook's default font metrics, so we
# pass in the character width that resolves to the closest achievable
# pixel width to the target.)
$overview.Range("E1").ColumnWidth = 16.25
$overview.Range("F1").ColumnWidth = 16.25
$zhcn.Range("C1").ColumnWidth = 16.25
$dede.Range("C1").ColumnWidth = 16.25
